$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new test row (row 11, test #9) describing the "8 gomb" Unit test.
$ws.Range("B11").Value = "Szicsák Bence"
$ws.Range("C11").Value = 44194
$ws.Range("D11").Value = "x"
$ws.Range("E11").Value = "Igen"
$ws.Range("F11").Value = "8 gomb"
$ws.Range("G11").Value = "Sikerült"

# "Sikerült" in the new row is highlighted in green, like a passed test.
$ws.Range("G11").Font.Name = "Times New Roman"
$ws.Range("G11").Font.Size = 12
$ws.Range("G11").Font.Color = 5287936

# Move the active selection to where the author left off editing.
[void]$ws.Range("I13").Select()
